$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 31, pushing the previous rows 31-35 down to 33-37.
$ws.Range("A31:A32").EntireRow.Insert()

# Row 31: new weekly entry - Angeleno / Especial
$ws.Range("A31").Value = 2
$ws.Range("B31").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C31").Value = 'Coquimbo'
$ws.Range("D31").Value = 44637
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 'Fruta'
$ws.Range("G31").Value = 100103
$ws.Range("H31").Value = 'Frutos de hueso (carozo)'
$ws.Range("I31").Value = 100103002
$ws.Range("J31").Value = 'Ciruela'
$ws.Range("K31").Value = 'Angeleno'
$ws.Range("L31").Value = 'Especial'
$ws.Range("M31").Value = 20
$ws.Range("N31").Value = 255000
$ws.Range("O31").Value = 260000
$ws.Range("P31").Value = 257500
$ws.Range("Q31").Value = '$/bins (450 kilos)'
$ws.Range("R31").Value = 'Región Metropolitana'
$ws.Range("S31").Value = 572
$ws.Range("T31").Value = 450

# Row 32: new weekly entry - Angeleno / Primera
$ws.Range("A32").Value = 2
$ws.Range("B32").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C32").Value = 'Coquimbo'
$ws.Range("D32").Value = 44637
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 'Fruta'
$ws.Range("G32").Value = 100103
$ws.Range("H32").Value = 'Frutos de hueso (carozo)'
$ws.Range("I32").Value = 100103002
$ws.Range("J32").Value = 'Ciruela'
$ws.Range("K32").Value = 'Angeleno'
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 20
$ws.Range("N32").Value = 225000
$ws.Range("O32").Value = 230000
$ws.Range("P32").Value = 227500
$ws.Range("Q32").Value = '$/bins (450 kilos)'
$ws.Range("R32").Value = 'Región Metropolitana'
$ws.Range("S32").Value = 506
$ws.Range("T32").Value = 450
